# Update column B (score) values on Sheet1 per "adding changes from last mtg".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B5"   = 5
    "B6"   = 4
    "B11"  = 5
    "B18"  = 2
    "B22"  = 1
    "B25"  = 5
    "B30"  = 1
    "B31"  = 5
    "B34"  = 5
    "B38"  = 4
    "B39"  = 4
    "B44"  = 5
    "B47"  = 2
    "B48"  = 4
    "B50"  = 2
    "B51"  = 2
    "B53"  = 2
    "B54"  = 2
    "B55"  = 5
    "B56"  = 4
    "B57"  = 4
    "B59"  = 2
    "B63"  = 5
    "B64"  = 5
    "B65"  = 4
    "B67"  = 4
    "B69"  = 5
    "B72"  = 2
    "B73"  = 4
    "B74"  = 4
    "B75"  = 5
    "B76"  = 5
    "B78"  = 2
    "B80"  = 5
    "B81"  = 2
    "B84"  = 5
    "B85"  = 5
    "B86"  = 5
    "B87"  = 4
    "B88"  = 4
    "B89"  = 5
    "B90"  = 5
    "B94"  = 2
    "B96"  = 2
    "B97"  = 5
    "B98"  = 2
    "B100" = 4
    "B101" = 5
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
